$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new daily row (row 76) under the existing D75 data block.
# Column A holds a date-looking string that must stay literal text (matching
# the source file's inlineStr cells), so force text formatting before the
# write and then drop back to the default "Normal" style so no stray
# number-format style is left attached to the cell.
$ws.Range("A76").NumberFormat = "@"
$ws.Range("A76").Value = "2025/10/08"
$ws.Range("A76").Style = "Normal"

$ws.Range("B76").Value = "水"
$ws.Range("C76").Value = 2
$ws.Range("D76").Value = 116
